$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2: Banana
$t.Cell(2, 2).Range.Text = "$ 0.54"
$t.Cell(2, 3).Range.Text = "1"
$t.Cell(2, 4).Range.Text = "$ 0.54"

# Row 3: Strawberry
$t.Cell(3, 2).Range.Text = "$ 0.28"
$t.Cell(3, 3).Range.Text = "4"
$t.Cell(3, 4).Range.Text = "$ 1.12"

# Row 4: Chicken
$t.Cell(4, 2).Range.Text = "$ 0.63"
$t.Cell(4, 3).Range.Text = "4"
$t.Cell(4, 4).Range.Text = "$ 2.52"

# Row 5: Bread
$t.Cell(5, 2).Range.Text = "$ 0.46"
$t.Cell(5, 3).Range.Text = "8"
$t.Cell(5, 4).Range.Text = "$ 3.68"

# Row 6: Eggs
$t.Cell(6, 2).Range.Text = "$ 0.01"
$t.Cell(6, 3).Range.Text = "6"
$t.Cell(6, 4).Range.Text = "$ 0.06"

# Row 7: Salad
$t.Cell(7, 2).Range.Text = "$ 0.37"
$t.Cell(7, 3).Range.Text = "7"
$t.Cell(7, 4).Range.Text = "$ 2.59"
